$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.177.77"
$ws.Range("E2").Value = "  -4.57%  "
$ws.Range("D3").Value = "1.654.14"
$ws.Range("E3").Value = "  -3.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.97"
$ws.Range("E5").Value = "  -4.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5137"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2595"
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06448"
$ws.Range("E9").Value = "  -3.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.96"
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07791"
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.296"
$ws.Range("E12").Value = "  -4.75%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.651.17"
$ws.Range("E13").Value = "  -3.84%  "
$ws.Range("D14").Value = "1.879.43"
$ws.Range("E14").Value = "  -3.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5528"
$ws.Range("E15").Value = "  -5.80%  "
$ws.Range("D16").Value = "0.0₅8022"
$ws.Range("E16").Value = "  -2.47%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.22"
$ws.Range("E17").Value = "  -5.77%  "
$ws.Range("D18").Value = "26.177.18"
$ws.Range("E18").Value = "  -4.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "210.71"
$ws.Range("E20").Value = "  -5.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.408"
$ws.Range("E21").Value = "  -5.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.08"
$ws.Range("E22").Value = "  -4.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.046"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.66"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.748"
$ws.Range("E26").Value = "  +3.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1177"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.980"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.86"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05099"
$ws.Range("E30").Value = "  -5.16%  "
$ws.Range("E31").Value = "  -4.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.361"
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.225"
$ws.Range("E33").Value = "  -6.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.563"
$ws.Range("E34").Value = "  -4.20%  "
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9250"
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5783"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("D39").Value = "1.164.83"
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01589"
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.562"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.663"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8240"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.27"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").Value = "1.790.07"
$ws.Range("E46").Value = "  -3.75%  "
$ws.Range("D47").Value = "0.0₈115"
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.53"
$ws.Range("E49").Value = "  -3.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.007"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.871"
$ws.Range("E51").Value = "  -2.96%  "
